$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("June")

# The day's log entry (row 2) shifts down to row 3 - relocate the whole
# row via Cut/paste so the "Annet" free-text note (which contains an
# embedded newline) is moved verbatim rather than rewritten.
$ws.Range("A2:K2").Cut($ws.Range("A3:K3"))

# Update the date label for the shifted entry (01.06.2018, Fri -> 02.06.2018, Sat)
# and flag the first three issue categories (Spindel, Spindel Oljetank,
# Hydraulisk Enhet) as having occurred that day.
$ws.Range("A3").Value = "02.06.2018, Sat"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
